$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "as at" date in the intro text (row 2)
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 15 November 2024"

# Remove the now-published entry (old row 5: "Mortgage and landlord possession statistics: July to
# September 2024", week commencing 11 Nov 2024). Everything below shifts up by one row.
$ws.Rows("5:5").Delete()

# The conditional formatting ranges referred to the old last row (80); after removing a row the
# table now ends at row 79, so re-point the two conditional formatting rule groups accordingly.
$cf = $ws.Cells.FormatConditions
$cf.Item(1).ModifyAppliesToRange($ws.Range("A5:F79"))
$cf.Item(4).ModifyAppliesToRange($ws.Range("A5:A79"))
